# Update "Horarios" workbook: refresh scrape timestamp, insert a new
# arrival row, and append a new arrival row on the two sheets that share
# the "14_ABASTO" / "215_ALUAR" dataset (LP1912 and 6203-6173).

$wb = $excel.ActiveWorkbook

$sheetNames = @("LP1912", "6203-6173")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Insert a brand-new row right above the current row 7 (00:46:06 | 01:58 | 14_ABASTO | 72),
    # shifting it (and the row below it) down by one.
    $ws.Rows.Item(7).Insert()

    # Fill in the newly inserted row 7.
    $ws.Cells.Item(7, 1).Value = "01:55:38"
    $ws.Cells.Item(7, 2).Value = "03:02"
    $ws.Cells.Item(7, 3).Value = "15_ABASTO"
    $ws.Cells.Item(7, 4).Value = 67

    # Append a new row after the current last data row (now row 9).
    $ws.Cells.Item(10, 1).Value = "01:55:38"
    $ws.Cells.Item(10, 2).Value = "03:48"
    $ws.Cells.Item(10, 3).Value = "14_ABASTO"
    $ws.Cells.Item(10, 4).Value = 113

    # Refresh the "last updated" stamp and row count.
    $ws.Range("A2").Value = "Última actualización: 01:55:38"
    $ws.Range("A3").Value = "Total filas: 5"
}
